$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.284.94"
$ws.Range("E2").Value = "  -4.05%  "
$ws.Range("D3").Value = "2.643.10"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "520.94"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "144.14"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "0.569"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("D9").Value = "6.67"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  -3.31%  "
$ws.Range("D11").Value = "0.338"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").Value = "3.111.16"
$ws.Range("E13").Value = "  -2.01%  "
$ws.Range("D14").Value = "58.270.95"
$ws.Range("E14").Value = "  -4.07%  "
$ws.Range("D15").Value = "20.82"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").Value = "2.648.91"
$ws.Range("E17").Value = "  -13.57%  "
$ws.Range("D18").Value = "338.16"
$ws.Range("E18").Value = "  -3.23%  "
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").Value = "10.45"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").Value = "6.29"
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "64.63"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("D24").Value = "0.424"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").Value = "0.167"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "0.0₃0797"
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").Value = "7.10"
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("D29").Value = "6.66"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").Value = "152.36"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").Value = "18.85"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("D34").Value = "4.13"
$ws.Range("E34").Value = "  -2.86%  "
$ws.Range("D35").Value = "0.909"
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("E36").Value = "  -4.65%  "
$ws.Range("D37").Value = "0.857"
$ws.Range("E37").Value = "  -2.90%  "
$ws.Range("D38").Value = "36.72"
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("E39").Value = "  -5.39%  "
$ws.Range("D40").Value = "3.64"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("D42").Value = "0.606"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("D44").Value = "269.98"
$ws.Range("E44").Value = "  -6.26%  "
$ws.Range("D45").Value = "19.42"
$ws.Range("E45").Value = "  -3.29%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "10.64"
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0536"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").Value = "2.036.83"
$ws.Range("E48").Value = "  -5.19%  "
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("E50").Value = "  -4.76%  "
$ws.Range("D51").Value = "18.33"
$ws.Range("E51").Value = "  -4.73%  "
